# Updated capital structure database
# Applies new figures to row 2 (aggregate, now count of 1 company),
# replaces row 3 with "Regional S.A.B. de C.V. (BMV:R A)" data,
# and removes the old standalone row 4 for that same company
# (its figures moved up into row 3, row count shrinks from 3 companies to 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"

$ws.Range("D2").Value = 0.144
$ws.Range("E2").Value = 0.145
$ws.Range("F2").Value = -0.007209999999999999
$ws.Range("K2").Value = 152.1
$ws.Range("L2").Value = 0.3661531054405393
$ws.Range("M2").Value = 3.1
$ws.Range("N2").Value = 0.002042429832652523
$ws.Range("O2").Value = 0.02038132807363577
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 3.1
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 344.5
$ws.Range("V2").Value = 0.2269732507576756
$ws.Range("W2").Value = 0.1597521268774288
$ws.Range("X2").Value = 0.128488543636297
$ws.Range("Y2").Value = 0.03126358324113179
$ws.Range("Z2").Value = 0.2473797046212482
$ws.Range("AB2").Value = 0.04985290552385511
$ws.Range("AC2").Value = -0.04985290552385511
$ws.Range("AD2").Value = 5446.2
$ws.Range("AF2").Value = 5446.2
$ws.Range("AG2").Value = 5101.7
$ws.Range("AH2").Value = 0.7820505456634118
$ws.Range("AI2").Value = 0.8464587121741969
$ws.Range("AJ2").Value = 0.7707077573834882
$ws.Range("AK2").Value = 0.8377725959012087

# --- Row 3 updates: Actinver data replaced with Regional S.A.B. de C.V. data ---
$ws.Range("B3").Value = "Regional S.A.B. de C.V. (BMV:R A)"

$ws.Range("D3").Value = 0.144
$ws.Range("E3").Value = 0.145
$ws.Range("F3").Value = -0.007209999999999999
$ws.Range("K3").Value = 152.1
$ws.Range("L3").Value = 0.3661531054405393
$ws.Range("M3").Value = 3.1
$ws.Range("N3").Value = 0.002042429832652523
$ws.Range("O3").Value = 0.02038132807363577
$ws.Range("S3").Value = 3.1
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 344.5
$ws.Range("V3").Value = 0.2269732507576756
$ws.Range("W3").Value = 0.1597521268774288
$ws.Range("X3").Value = 0.128488543636297
$ws.Range("Y3").Value = 0.03126358324113179
$ws.Range("Z3").Value = 0.2473797046212482
$ws.Range("AB3").Value = 0.04985290552385511
$ws.Range("AC3").Value = -0.04985290552385511
$ws.Range("AD3").Value = 5446.2
$ws.Range("AF3").Value = 5446.2
$ws.Range("AG3").Value = 5101.7
$ws.Range("AH3").Value = 0.7820505456634118
$ws.Range("AI3").Value = 0.8464587121741969
$ws.Range("AJ3").Value = 0.7707077573834882
$ws.Range("AK3").Value = 0.8377725959012087

# --- Remove the now-duplicate row 4 (Regional S.A.B. de C.V.) ---
$ws.Rows.Item(4).Delete()
